$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 14 (TC_013): the test result is reverted from a passed/verified
# outcome to a "not executed" outcome.
$ws.Range("H14").Value = "Test not executed"
$ws.Range("I14").Value = "Not Run"

# I14 previously carried the green "PASSED" highlight style (s=4). Copy
# the plain formatting already used by H14 (s=3, no fill) onto I14 so it
# loses the green highlight and matches the surrounding cell style.
$ws.Range("H14").Copy()
$ws.Range("I14").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 28 (TC_027): fix quoting/casing in the test objective text.
$ws.Range("C28").Value = "Verify visibility of ""Don't have an account?"" text on login page."
